$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column Y ("from_scratch") ---------------------------------
# Header cell Y1: copy header formatting from X1, then set its text.
$ws.Range("X1").Copy()
$ws.Range("Y1").PasteSpecial(-4122)
$ws.Range("Y1").Value = "from_scratch"

# Data cells Y2:Y31 on the existing rows are written as blank/empty text
# (mirrors the existing empty W/X inlineStr cells in this sheet). Typing a
# lone apostrophe enters an empty, left-as-text value; resetting the style
# back to Normal drops the "quote prefix" formatting that the apostrophe
# would otherwise leave behind.
for ($r = 2; $r -le 31; $r++) {
  $cell = $ws.Cells.Item($r, 25)
  $cell.Value = "'"
  $cell.Style = "Normal"
}

# --- X31 switches from a boolean TRUE to the number 1 -------------------
$ws.Range("X31").Value = 1

# --- Append three new rows of run data (32-34) --------------------------
$newRows = @(
  @($false, "runs_report.xlsx", "['050']", 500, 0.9, 3, "uniform", "linear", "mlp_lin", "mlp", "fixed-grid", "rk4", "mlp", "[1]", "tanh", "none", "[1, 3, 1]", 10,  3000, 0.001, "run_031", "xavier", 0, 1),
  @($false, "runs_report.xlsx", "['050']", 500, 0.9, 3, "uniform", "linear", "mlp_lin", "mlp", "fixed-grid", "rk4", "mlp", "[1]", "tanh", "none", "[1, 3, 1]", 10,  3000, 0.001, "run_032", "xavier", 0, 0),
  @($false, "runs_report.xlsx", "['050']", 500, 0.9, 3, "uniform", "linear", "mlp_lin", "mlp", "fixed-grid", "rk4", "mlp", "[1]", "tanh", "none", "[1, 3, 1]", 200, 3000, 0.001, "run_033", "xavier", 0, 0)
)

$row = 32
foreach ($data in $newRows) {
  $ws.Cells.Item($row, 1).Value  = $data[0]
  $ws.Cells.Item($row, 2).Value  = $data[1]

  # Column C ("retrain") is left blank/empty text on these new rows.
  $cCell = $ws.Cells.Item($row, 3)
  $cCell.Value = "'"
  $cCell.Style = "Normal"

  $ws.Cells.Item($row, 4).Value  = $data[2]
  $ws.Cells.Item($row, 5).Value  = $data[3]
  $ws.Cells.Item($row, 6).Value  = $data[4]
  $ws.Cells.Item($row, 7).Value  = $data[5]
  $ws.Cells.Item($row, 8).Value  = $data[6]
  $ws.Cells.Item($row, 9).Value  = $data[7]
  $ws.Cells.Item($row, 10).Value = $data[8]
  $ws.Cells.Item($row, 11).Value = $data[9]
  $ws.Cells.Item($row, 12).Value = $data[10]
  $ws.Cells.Item($row, 13).Value = $data[11]
  $ws.Cells.Item($row, 14).Value = $data[12]
  $ws.Cells.Item($row, 15).Value = $data[13]
  $ws.Cells.Item($row, 16).Value = $data[14]
  $ws.Cells.Item($row, 17).Value = $data[15]
  $ws.Cells.Item($row, 18).Value = $data[16]
  $ws.Cells.Item($row, 19).Value = $data[17]
  $ws.Cells.Item($row, 20).Value = $data[18]
  $ws.Cells.Item($row, 21).Value = $data[19]
  $ws.Cells.Item($row, 22).Value = $data[20]
  $ws.Cells.Item($row, 23).Value = $data[21]
  $ws.Cells.Item($row, 24).Value = $data[22]
  $ws.Cells.Item($row, 25).Value = $data[23]

  $row++
}
